$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'67.847.47"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'  +1.00%  "
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'2.626.74"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'  +0.42%  "
$ws.Range("E3").Style = "Normal"
$ws.Range("E4").Value = "'  -0.05%  "
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'597.90"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'  +0.42%  "
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = "'153.68"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'  +0.77%  "
$ws.Range("E6").Style = "Normal"
$ws.Range("E7").Value = "'  -0.05%  "
$ws.Range("E7").Style = "Normal"
$ws.Range("E8").Value = "'  -1.41%  "
$ws.Range("E8").Style = "Normal"
$ws.Range("D9").Value = "'2.626.93"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'  +0.51%  "
$ws.Range("E9").Style = "Normal"
$ws.Range("D10").Value = "'0.134"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'  +10.00%  "
$ws.Range("E10").Style = "Normal"
$ws.Range("E11").Value = "'  -0.69%  "
$ws.Range("E11").Style = "Normal"
$ws.Range("E12").Value = "'  +0.77%  "
$ws.Range("E12").Style = "Normal"
$ws.Range("E13").Value = "'  +0.21%  "
$ws.Range("E13").Style = "Normal"
$ws.Range("D14").Value = "'27.53"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'  +0.11%  "
$ws.Range("E14").Style = "Normal"
$ws.Range("E15").Value = "'  +4.61%  "
$ws.Range("E15").Style = "Normal"
$ws.Range("D16").Value = "'3.103.81"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'  +0.29%  "
$ws.Range("E16").Style = "Normal"
$ws.Range("D17").Value = "'67.783.11"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'  +1.16%  "
$ws.Range("E17").Style = "Normal"
$ws.Range("D18").Value = "'2.624.29"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'  +0.35%  "
$ws.Range("E18").Style = "Normal"
$ws.Range("D19").Value = "'11.44"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'  +3.75%  "
$ws.Range("E19").Style = "Normal"
$ws.Range("D20").Value = "'370.81"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'  +1.97%  "
$ws.Range("E20").Style = "Normal"
$ws.Range("D21").Value = "'7.43"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'  +0.20%  "
$ws.Range("E21").Style = "Normal"
$ws.Range("D22").Value = "'4.23"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'  -1.31%  "
$ws.Range("E22").Style = "Normal"
$ws.Range("E23").Value = "'  -1.45%  "
$ws.Range("E23").Style = "Normal"
$ws.Range("D24").Value = "'2.06"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'  -0.30%  "
$ws.Range("E24").Style = "Normal"
$ws.Range("D25").Value = "'72.20"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'  +1.41%  "
$ws.Range("E25").Style = "Normal"
$ws.Range("E26").Value = "'  +0.08%  "
$ws.Range("E26").Style = "Normal"
$ws.Range("D27").Value = "'9.87"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "'  -0.46%  "
$ws.Range("E27").Style = "Normal"
$ws.Range("D28").Value = "'2.754.36"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "'  +0.40%  "
$ws.Range("E28").Style = "Normal"
$ws.Range("E29").Value = "'  +2.80%  "
$ws.Range("E29").Style = "Normal"
$ws.Range("E30").Value = "'  +0.19%  "
$ws.Range("E30").Style = "Normal"
$ws.Range("D31").Value = "'573.54"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "'  -1.19%  "
$ws.Range("E31").Style = "Normal"
$ws.Range("D32").Value = "'7.89"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "'  +1.47%  "
$ws.Range("E32").Style = "Normal"
$ws.Range("D33").Value = "'1.40"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "'  +0.63%  "
$ws.Range("E33").Style = "Normal"
$ws.Range("E34").Value = "'  +1.18%  "
$ws.Range("E34").Style = "Normal"
$ws.Range("E35").Value = "'  -0.06%  "
$ws.Range("E35").Style = "Normal"
$ws.Range("D36").Value = "'0.126"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "'  +1.68%  "
$ws.Range("E36").Style = "Normal"
$ws.Range("D37").Value = "'1.52"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "'  +1.71%  "
$ws.Range("E37").Style = "Normal"
$ws.Range("D38").Value = "'159.46"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "'  +1.33%  "
$ws.Range("E38").Style = "Normal"
$ws.Range("D39").Value = "'19.14"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "'  +0.68%  "
$ws.Range("E39").Style = "Normal"
$ws.Range("E40").Value = "'  +4.89%  "
$ws.Range("E40").Style = "Normal"
$ws.Range("E41").Value = "'  +0.35%  "
$ws.Range("E41").Style = "Normal"
$ws.Range("D42").Value = "'5.34"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'  +1.93%  "
$ws.Range("E42").Style = "Normal"
$ws.Range("D43").Value = "'0.0₆0335"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'  +17.33%  "
$ws.Range("E43").Style = "Normal"
$ws.Range("D44").Value = "'2.62"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'  +3.16%  "
$ws.Range("E44").Style = "Normal"
$ws.Range("D45").Value = "'17.38"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'  +6.19%  "
$ws.Range("E45").Style = "Normal"
$ws.Range("E46").Value = "'  +0.07%  "
$ws.Range("E46").Style = "Normal"
$ws.Range("D47").Value = "'40.25"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "'  -2.15%  "
$ws.Range("E47").Style = "Normal"
$ws.Range("D48").Value = "'155.55"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'  -0.66%  "
$ws.Range("E48").Style = "Normal"
$ws.Range("E49").Value = "'  -0.81%  "
$ws.Range("E49").Style = "Normal"
$ws.Range("E50").Value = "'  +0.84%  "
$ws.Range("E50").Style = "Normal"
$ws.Range("E51").Value = "'  -0.11%  "
$ws.Range("E51").Style = "Normal"
